$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain text (not numeric-looking), set directly.
$plainUpdates = @{
    "D2" = "60.191.68"
    "E2" = "  -0.60%  "
    "D3" = "2.620.60"
    "E3" = "  +0.98%  "
    "E4" = "  -0.09%  "
    "E5" = "  +1.06%  "
    "E6" = "  -2.88%  "
    "E7" = "  +0.02%  "
    "E8" = "  -4.17%  "
    "D9" = "2.622.85"
    "E9" = "  +0.70%  "
    "E10" = "  -5.38%  "
    "E11" = "  +0.49%  "
    "E12" = "  -1.38%  "
    "E13" = "  -0.70%  "
    "D14" = "3.073.86"
    "E14" = "  +0.76%  "
    "D15" = "60.226.35"
    "E15" = "  -0.51%  "
    "E16" = "  -2.20%  "
    "E17" = "  -1.49%  "
    "D18" = "2.614.82"
    "E18" = "  +0.53%  "
    "E19" = "  -2.24%  "
    "E20" = "  -2.64%  "
    "E21" = "  -1.33%  "
    "E22" = "  -1.71%  "
    "E23" = "  -0.17%  "
    "E24" = "  -0.81%  "
    "E25" = "  -1.93%  "
    "E26" = "  +0.00%  "
    "E27" = "  -2.10%  "
    "D28" = "0.0₃0811"
    "E28" = "  -3.61%  "
    "E29" = "  -3.57%  "
    "E30" = "  +0.00%  "
    "E31" = "  -4.16%  "
    "E32" = "  -0.09%  "
    "E33" = "  -2.26%  "
    "E34" = "  -0.65%  "
    "E35" = "  -3.44%  "
    "E36" = "  +0.13%  "
    "E37" = "  -4.70%  "
    "E38" = "  +3.14%  "
    "E40" = "  -3.17%  "
    "E41" = "  -3.92%  "
    "E42" = "  +1.07%  "
    "E43" = "  +0.05%  "
    "E44" = "  -1.06%  "
    "E45" = "  -0.03%  "
    "E46" = "  -2.03%  "
    "E47" = "  -0.44%  "
    "E48" = "  +0.89%  "
    "E49" = "  -2.02%  "
    "E50" = "  -2.07%  "
    "D51" = "1.960.61"
    "E51" = "  -0.49%  "
}

# Cells whose new value looks like a plain number (e.g. "522.28") but must
# remain a text string (matching the original inline-string cell type).
# Temporarily format as Text so Excel does not coerce the value to a number,
# then restore the original (default) cell style.
$textForcedUpdates = @{
    "D5" = "522.28"
    "D6" = "148.96"
    "D7" = "0.999"
    "D8" = "0.571"
    "D10" = "6.31"
    "D16" = "21.22"
    "D20" = "342.48"
    "D21" = "10.44"
    "D22" = "6.12"
    "D23" = "0.998"
    "D24" = "60.56"
    "D26" = "1.00"
    "D31" = "6.01"
    "D33" = "18.98"
    "D34" = "149.68"
    "D35" = "3.97"
    "D36" = "0.917"
    "D38" = "0.863"
    "D39" = "36.45"
    "D41" = "3.63"
    "D42" = "289.91"
    "D43" = "0.624"
    "D44" = "0.101"
    "D45" = "0.998"
    "D46" = "0.0548"
    "D47" = "19.43"
    "D50" = "4.72"
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

foreach ($cellRef in $textForcedUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$cellRef]
    $cell.Style = "Normal"
}
